# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '61.999.33'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E2').Value = '  +5.93%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '3.072.90'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('E4').Value = '  +0.22%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '579.04'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E5').Value = '  +3.69%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '143.55'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E6').Value = '  +5.41%  '
$ws.Range('E7').Value = '  +0.18%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '3.057.69'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E8').Value = '  +2.87%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.524'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E9').Value = '  +1.84%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.139'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E10').Value = '  +6.75%  '
$ws.Range('E11').Value = '  +13.68%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.466'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E12').Value = '  +2.80%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.0000243'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E13').Value = '  +7.80%  '
$ws.Range('E14').Value = '  +4.10%  '
$ws.Range('E15').Value = '  -0.19%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '3.571.73'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E16').Value = '  +4.10%  '
$ws.Range('E17').Value = '  +4.14%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '3.062.54'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E18').Value = '  +3.14%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '61.720.77'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E19').Value = '  +5.37%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '450.00'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E20').Value = '  +6.97%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '13.90'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('E22').Value = '  +3.24%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '7.31'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E23').Value = '  +3.18%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '13.66'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E24').Value = '  +2.87%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '82.07'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E25').Value = '  +2.71%  '
$ws.Range('E26').Value = '  -0.03%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E27').Value = '  +7.10%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E28').Value = '  -0.11%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.64'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E29').Value = '  +5.40%  '
$ws.Range('E30').Value = '  +4.72%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '6.52'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E31').Value = '  +7.92%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '26.59'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E32').Value = '  +3.91%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.107'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E33').Value = '  +8.38%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.0₃0824'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E34').Value = '  +9.76%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.03'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('E36').Value = '  +6.09%  '
$ws.Range('E37').Value = '  +6.44%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '50.20'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E38').Value = '  +3.30%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.99'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E39').Value = '  +10.17%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '8.84'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E40').Value = '  +1.82%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '418.07'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E41').Value = '  +5.40%  '
$ws.Range('E42').Value = '  +5.48%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.781.91'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E43').Value = '  +1.65%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('E45').Value = '  +9.21%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '37.45'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E46').Value = '  +17.08%  '
$ws.Range('E47').Value = '  +5.67%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '123.75'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('E50').Value = '  +2.05%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '24.32'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range('E51').Value = '  +5.14%  '
